# Add a "target_choice" column (F) to the bias parameters table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new column.
$ws.Range("F1").Value = "target_choice"

# Values for the new column, one per data row.
$ws.Range("F2").Value = "A"
$ws.Range("F3").Value = "B"
$ws.Range("F4").Value = "A"

# The original A2 cell carried a stray numeric-format style (xfId 1) that
# is not actually applied (format stays "General"); drop it so A2 goes
# back to the default style, matching the rest of the column.
$ws.Range("A2").ClearFormats()

# Move the selection to the new last cell, as in the edited workbook.
$ws.Range("F4").Select()
